$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row above row 12 (shifts rows 12:24 down to 13:25,
# copying formatting from the row above, matching Excel's native
# "Insert Row" behaviour).
$ws.Rows.Item(12).Insert()

# Populate the new row with the added translation pair
# (checkbox label shown/hidden on the task bar).
$ws.Range("C12").Value = "Toon niet op Taakbalk"
$ws.Range("D12").Value = "Hide from task bar"

# Move the active selection to the newly added cell D12.
$ws.Range("D12").Select()
